$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.085905909538269
$ws.Range("B1").Value = 2.149283170700073
$ws.Range("C1").Value = 9.398229598999023
$ws.Range("D1").Value = 1.026668548583984
$ws.Range("E1").Value = 1.088949918746948
